$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New time-tracking entry in row 43 ---
# A43: date 2018-10-16 (serial 43389), keeps existing date-format style
$ws.Range("A43").Value = 43389
# B43: worker name
$ws.Range("B43").Value = "Daniel"
# C43: description of the work done
$ws.Range("C43").Value = "Mit Demo von Hr. Prof. Stütz Fehler gefixed"
# D43: hours spent
$ws.Range("D43").Value = 1.5

# H3 contains =SUMIF(B2:B150,"Daniel",D2:D150) and recalculates automatically
# once the new row's data is present (62 -> 63.5).

# --- Update the view's scroll position / selection ---
$ws.Range("E11").Select() | Out-Null
